# Apply the cryptos price-list refresh described in the commit:
#   "Updated cryptos list on Sat May 27 18:40:28 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Price" column (D) -------------------------------------------------
# Several of the new price strings look like plain decimal numbers
# (e.g. "1.015", "309.80", "5.400"). Assigning such a string straight to
# .Value lets Excel auto-convert it to a floating point number (losing the
# trailing zero / original text). To keep these as literal text - exactly
# like the inline strings already in the sheet - we build the text with a
# text formula and then paste back only the resulting value, which keeps
# the cell as a text cell without touching its number format/style.
$ws.Range("D2").Formula = "=""27.168.43"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("D3").Formula = "=""1.852.41"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("D4").Formula = "=""1.015"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""309.80"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D7").Formula = "=""0.4776"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("D8").Formula = "=""0.3692"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("D9").Formula = "=""0.07257"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D10").Formula = "=""0.9329"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""19.91"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D12").Formula = "=""0.07797"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("D13").Formula = "=""1.830.45"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("D14").Formula = "=""5.400"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D15").Formula = "=""6.488"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("D16").Formula = "=""89.43"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("D17").Formula = "=""1.017"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("D18").Formula = "=""0.000008705"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=""27.163.79"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D21").Formula = "=""14.62"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("D22").Formula = "=""5.063"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D23").Formula = "=""10.66"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D24").Formula = "=""1.947"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=""153.11"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D26").Formula = "=""18.36"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""1.989"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""114.83"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""4.938"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D30").Formula = "=""0.08883"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D31").Formula = "=""3.299"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""1.184"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""4.527"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("D34").Formula = "=""0.7406"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("D35").Formula = "=""2.683"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("D36").Formula = "=""1.117"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D37").Formula = "=""0.01985"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("D38").Formula = "=""0.05278"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("D39").Formula = "=""2.977"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=""0.5291"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D41").Formula = "=""7.048"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=""0.1526"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D43").Formula = "=""8.313"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("D44").Formula = "=""10.57"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("D45").Formula = "=""0.4745"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("D46").Formula = "=""1.015"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("D47").Formula = "=""102.16"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("D48").Formula = "=""1.621"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("D49").Formula = "=""65.89"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("D50").Formula = "=""0.06065"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("D51").Formula = "=""0.8936"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

# --- "Volume(1h)" column (E) ---------------------------------------------
# These values are always wrapped in spaces and end with a "%" sign, so
# they can never be misread as numbers - a direct text assignment is safe.
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("E7").Value = "  +2.27%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +1.94%  "
$ws.Range("E12").Value = "  +1.25%  "
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  -3.48%  "
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("E37").Value = "  +2.20%  "
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("E51").Value = "  +0.79%  "

$excel.CutCopyMode = 0
